# Catalogue.xlsx edit:
# 1. "Variabele componenten in excel kunnen inlezen" -> correct the CPU
#    sheet's integrated-graphics description so it is brand/model neutral
#    ("Intel Graphics 4600, 4 cores" -> "iGPU 4600, 4 cores") for both CPUs
#    that shipped with that iGPU (i5 4670K, i7 4770K).
# 2. "Cleaning" -> tidy up the workbook/cell style naming and leave the
#    CPU sheet active (instead of GPU) with the selection parked past the
#    data range, matching the author's last view state.

$wb = $excel.ActiveWorkbook

# --- 1. Fix the CPU integrated-graphics description -----------------------
$cpu = $wb.Worksheets.Item("CPU")
$cpu.Range("D2").Value = "iGPU 4600, 4 cores"
$cpu.Range("D3").Value = "iGPU 4600, 4 cores"

# --- 2. Cleaning ------------------------------------------------------------
# Rename the default "Normal" cell style to the Dutch-locale "Standaard".
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Name = "Standaard"

# Move the active tab / tab-selected state from GPU to CPU, and leave the
# selection on CPU!D12 (GPU regains its plain, unselected sheetView).
$cpu.Activate()
$cpu.Range("D12").Select()
